$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - new daily worklog entry
$ws.Cells.Item(3,1).Value = "06f61e43-6fc6-4ad2-8b18-68290b7dc68c"
$ws.Cells.Item(3,2).Value = "b4ee870e-ee7a-4152-b5d3-f40ee4f6bfb6"
$ws.Cells.Item(3,3).Value = $ws.Cells.Item(2,3).Value2
$ws.Cells.Item(2,4).Copy($ws.Cells.Item(3,4))
$ws.Cells.Item(3,4).Value = 45469.50799246528
$ws.Cells.Item(2,5).Copy($ws.Cells.Item(3,5))
$ws.Cells.Item(3,5).Value = 45469.50807127314
$ws.Cells.Item(3,6).Value = 0
$ws.Cells.Item(3,7).Value = $ws.Cells.Item(2,7).Value2
$ws.Cells.Item(3,8).Value = $ws.Cells.Item(2,8).Value2

# Row 4 - new daily worklog entry
$ws.Cells.Item(4,1).Value = "81646156-cfb4-4182-a5d0-8f76f3b681fe"
$ws.Cells.Item(4,2).Value = "b4ee870e-ee7a-4152-b5d3-f40ee4f6bfb6"
$ws.Cells.Item(4,3).Value = $ws.Cells.Item(2,3).Value2
$ws.Cells.Item(2,4).Copy($ws.Cells.Item(4,4))
$ws.Cells.Item(4,4).Value = 45469.50818040509
$ws.Cells.Item(2,5).Copy($ws.Cells.Item(4,5))
$ws.Cells.Item(4,5).Value = 45469.508230810185
$ws.Cells.Item(4,6).Value = 0
$ws.Cells.Item(4,7).Value = $ws.Cells.Item(2,7).Value2
$ws.Cells.Item(4,8).Value = $ws.Cells.Item(2,8).Value2
